$wb = $excel.ActiveWorkbook

# Locate the most recent existing weekly sheet (2025-11-05) to anchor the
# new sheet after it, and to source the header formatting from.
$srcSheet = $wb.Worksheets.Item(3)

# --- 1) Add the new weekly sheet right after it ---------------------------
$newSheetName = "2025-11-12"
$newWs = $wb.Worksheets.Add($null, $srcSheet)
$newWs.Name = $newSheetName

# --- 2) Header row ----------------------------------------------------------
$newWs.Cells.Item(1, 1).Value = "rank"
$newWs.Cells.Item(1, 2).Value = "title"
$newWs.Cells.Item(1, 3).Value = "volume"
$newWs.Cells.Item(1, 4).Value = "publisher"

# Reuse the exact header formatting (bold font + border + center/top align)
# already used on the other weekly sheets instead of re-deriving it.
$srcSheet.Range("A1:D1").Copy()
$newWs.Range("A1:D1").PasteSpecial(-4122)

# --- 3) Ranking rows --------------------------------------------------------
# Each entry: title, volume
$data = @(
  @("ワンパンマン", 35),
  @("BLUE GIANT MOMENTUM", 6),
  @("九条の大罪", 15),
  @("少年院ウシジマくん", 6),
  @("傷モノの花嫁", 9),
  @("転生したらスライムだった件", 30),
  @("SPY×FAMILY", 16),
  @("転生したらスライムだった件~魔物の国の歩き方~", 8),
  @("永年雇用は可能でしょうか", 5),
  @("エルフ先生と呼ばないで! 第1話", 1),
  @("ギャラ飲み女子とラーメンおじさん", 1),
  @("逆行した元悪役令嬢、性格の悪さは直さず処刑エンド回避します! 第1話(アリアンローズコミックス)", 1),
  @("異世界整体師 ~美女も亜人も魔物も竜も、お前ら全員揉みほぐす!!~", 1),
  @("終わりのセラフ", 35),
  @("無職転生 ~異世界行ったら本気だす~ 失意の魔術師編", 1),
  @("聖女様になりたいのに攻撃魔法しか使えないんですけど!? 第2話", 2),
  @("まほまね", 1),
  @("熱造カノジョ1", 1),
  @("ガチャを回して仲間を増やす 最強の美少女軍団を作り上げろ THE COMIC", 3),
  @("賭ケグルイ", 20),
  @("勇者は魔王が好きらしい", 1),
  @("ハイブルク家三男は小悪魔ショタです 第1話", 1),
  @("勇者に殺される悪徳領主に転生した俺、序盤に鍛えすぎたせいで勇者の地位を奪ってしまう1", 1),
  @("逆立てないで!カルマくん1", 1),
  @("異界から聖女が来たのでお役御免になりました~処刑されそうなので隠した力を解放させていただきます!~1", 1),
  @("このマンガがすごい! comics おとなしく泣き寝入りするとでも思いましたか?第1話①", 1),
  @("スキ妻とモチ夫~うちの妻が愛らしすぎて心配です!~", 1),
  @("「追放村」領主の超開拓 ~追放者だらけの辺境村がやがて世界に覇権を唱えるようです~(コミック) 1話", 1),
  @("鵺の陰陽師", 12),
  @("世界の終わりのいずこねこ 完全版 茉里ver.", 1),
  @("ディエンビエンフー・プレス 完全版", 1),
  @("世界の終わりの魔法使い 完全版 1 すべての始まり", 1),
  @("アオザイ通信 完全版 1 食と文化", 1),
  @("ヤング・アライブ・イン・ラブ 完全版1", 1),
  @("時ノ檻~死に戻りの霊装使い、能力で最強へと至る~1", 1),
  @("お茶屋さんは賢者見習い", 1),
  @("夜叉王の最愛 ~虐げられた治癒の乙女は溺愛される~", 1),
  @("宮中は噂のたえない職場にて", 1),
  @("RED&BLUE 第5話", 5),
  @("復讐完遂者の人生二周目異世界譚 THE COMIC", 2),
  @("あやかし妖怪娘", 1),
  @("ギルドの受付嬢は定時上がりの夢を見る@COMIC 第1話", 1),
  @("摩天楼", 1),
  @("転生して田舎でスローライフをおくりたい", 14),
  @("オーイ! とんぼ", 59),
  @("ワンパンマン", 34),
  @("片田舎のおっさん、剣聖になる~ただの田舎の剣術師範だったのに、大成した弟子たちが俺を放ってくれない件~(話売り) #40", 40),
  @("祝福のチェスカ : 12 エピローグ", 14),
  @("魔導具師ダリヤはうつむかない ~王立高等学院編~", 1),
  @("異世界帰りのパラディンは、最強の除霊師となる", 1),
  @("私が恋したきみじゃない", 2),
  @("「お前ごときが魔王に勝てると思うな」と勇者パーティを追放されたので、王都で気ままに暮らしたい THE COMIC", 1),
  @("毎日家に来るギャルが距離感ゼロでも優しくない THE COMIC", 1),
  @("婚約破棄されたら異国の王子に溺愛されました 甘~いキスは悦楽の予感", 1),
  @("このマンガがすごい! comics おとなしく泣き寝入りするとでも思いましたか?第1話②", 2),
  @("このマンガがすごい! comics おとなしく泣き寝入りするとでも思いましたか?第1話③", 3),
  @("ぼくは君の奇跡 ~The Miracle of Teddy Bear~", 1),
  @("ディエンビエンフー 完全版", 1),
  @("西島大介短編集 1 夏の彗星", 1),
  @("竜騎士のお気に入り:", 10),
  @("月が導く異世界道中16", 16),
  @("MAJOR 2nd(メジャーセカンド)", 31),
  @("ルリドラゴン", 4),
  @("ディエンビエンフー・プレス 完全版", 2),
  @("世界の終わりの魔法使い 完全版 2 恋におちた悪魔", 2),
  @("世界の終わりの魔法使い 完全版 3 影の子どもたち", 3),
  @("アオザイ通信 完全版 2 歴史と戦争", 2),
  @("若き社長は婚約者の姉を溺愛する1", 1),
  @("フリースキルで最強冒険者 ~ペットも無双で異世界生活が楽しすぎる~", 1),
  @("治癒魔法の間違った使い方 ~誘いの街・レストバレー~", 1),
  @("闇金クロサキ~復讐は計画的に~1", 1),
  @("イリーガル 外道には制裁を1", 1),
  @("私以外みんな幸せ", 1),
  @("東京婚活難民 ~結婚なんて、その気になればすぐできる…と思ってた~", 1),
  @("賢者の弟子を名乗る賢者 THE COMIC", 1),
  @("元最強の剣士は、異世界魔法に憧れる THE COMIC", 1),
  @("ラッキードッグ", 1),
  @("経験人数が見えるメガネ", 1),
  @("催芽さんは催眠アプられたい1", 1),
  @("暴食のベルセルク~俺だけレベルという概念を突破する~ THE COMIC", 3),
  @("異世界で聖女になったので国民的アイドルを目指します", 1),
  @("異界から聖女が来たのでお役御免になりました~処刑されそうなので隠した力を解放させていただきます!~2", 2),
  @("セカンドショジョと契約彼氏", 6),
  @("玉の輿ゲーム", 4),
  @("玉の輿ゲーム", 5),
  @("ぼくは君の奇跡 ~The Miracle of Teddy Bear~", 2),
  @("ぼくは君の奇跡 ~The Miracle of Teddy Bear~", 3),
  @("ぼくは君の奇跡 ~The Miracle of Teddy Bear~", 4),
  @("シャングリラ・フロンティア ~クソゲーハンター、神ゲーに挑まんとす~", 24),
  @("追放されるたびにスキルを手に入れた俺が、100の異世界で2周目無双", 5),
  @("神々の加護で生産革命~異世界の片隅でまったりスローライフしてたら、なぜか多彩な人材が集まって最強国家ができてました~(コミック)", 6),
  @("ありす、宇宙までも", 5),
  @("今どきの若いモンは", 29),
  @("カグラバチ", 9),
  @("悪祓士のキヨシくん", 6),
  @("あかね噺", 19),
  @("社内探偵", 65),
  @("俺の声に堕ちてください", 14),
  @("死にかけ令嬢ですが冷徹な騎士様と理想の殿方を探します!1", 1),
  @("じゃじゃ馬令嬢の婚活は前途多難です~辺境伯の筆頭護衛を攻略できません!~1", 1)
)

$highlightColor = 13499135  # RGB(255,250,205) -> "FFFACD", light-yellow new-entry highlight

$row = 2
foreach ($entry in $data) {
    $title = $entry[0]
    $volume = $entry[1]

    $newWs.Cells.Item($row, 1).Value = $row - 1
    $newWs.Cells.Item($row, 2).Value = $title
    $newWs.Cells.Item($row, 3).Value = $volume

    if ($volume -le 3) {
        $newWs.Cells.Item($row, 3).Interior.Color = $highlightColor
    }

    $row = $row + 1
}

# --- 4) Drop the now-unused blank "publisher" placeholder cells on the ------
#        previous week's sheet (2025-11-05) -- they were inline-string stubs
#        with no actual content.
$srcSheet.Range("D2:D101").ClearContents()

Write-Host "Added sheet '$newSheetName' with $($row - 2) ranking rows."
